$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("For plotting")

$ws.Range("D2").Value = 0.266347322993771
$ws.Range("E2").Value = 0.382308228705411

$ws.Range("D3").Value = 0.171051162687793
$ws.Range("E3").Value = 0.339831382016731

$ws.Range("D4").Value = 0.390697967538696
$ws.Range("E4").Value = 0.494901435241087
